$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create row 43 (new row), copying formatting from row 42 ---
$ws.Range("A42:E42").Copy() | Out-Null
$ws.Range("A43:E43").PasteSpecial(-4122) | Out-Null

# --- Update merged cells: remove stale merges ---
$ws.Range("A2:A21").UnMerge()
$ws.Range("A22:A42").UnMerge()
$ws.Range("B10:B15").UnMerge()
$ws.Range("B16:B17").UnMerge()
$ws.Range("B18:B19").UnMerge()
$ws.Range("B24:B25").UnMerge()
$ws.Range("B27:B28").UnMerge()
$ws.Range("B29:B30").UnMerge()
$ws.Range("B31:B32").UnMerge()
$ws.Range("B34:B39").UnMerge()
$ws.Range("C10:C14").UnMerge()

# --- Update merged cells: add new merges ---
$ws.Range("A2:A20").Merge()
$ws.Range("A21:A43").Merge()
$ws.Range("B9:B14").Merge()
$ws.Range("B17:B18").Merge()
$ws.Range("B21:B22").Merge()
$ws.Range("B25:B26").Merge()
$ws.Range("B28:B29").Merge()
$ws.Range("B30:B31").Merge()
$ws.Range("B32:B33").Merge()
$ws.Range("B35:B40").Merge()
$ws.Range("C9:C13").Merge()

# --- Update cell values to match target state ---
$ws.Range("E2").Value = 18350
$ws.Range("E3").Value = 46050
$ws.Range("E4").Value = 9610
$ws.Range("E5").Value = -225
$ws.Range("E6").Value = 1470
$ws.Range("B7").Value = "Other Current Assets"
$ws.Range("C7").Value = "Miscellaneous Current Assets"
$ws.Range("E7").Value = 3410
$ws.Range("B8").Value = "Total Current Assets"
$ws.Range("C8").Value = ""
$ws.Range("E8").Value = 78660
$ws.Range("B9").Value = "Net Property, Plant & Equipment"
$ws.Range("C9").Value = "Property, Plant & Equipment - Gross"
$ws.Range("D9").Value = "Buildings"
$ws.Range("E9").Value = 13330
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = "Construction in Progress"
$ws.Range("E10").Value = ""
$ws.Range("D11").Value = "Leases"
$ws.Range("D12").Value = "Computer Software and Equipment"
$ws.Range("E12").Value = 10920
$ws.Range("D13").Value = "Other Property, Plant & Equipment"
$ws.Range("E13").Value = 1950
$ws.Range("C14").Value = "Accumulated Depreciation"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = 8860
$ws.Range("B15").Value = "Total Investments and Advances"
$ws.Range("C15").Value = "LT Investment - Affiliate Companies"
$ws.Range("E15").Value = ""
$ws.Range("B16").Value = "Long-Term Note Receivable"
$ws.Range("C16").Value = ""
$ws.Range("E16").Value = 1330
$ws.Range("B17").Value = "Intangible Assets"
$ws.Range("C17").Value = "Net Goodwill"
$ws.Range("E17").Value = 15600
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = "Net Other Intangibles"
$ws.Range("E18").Value = 4610
$ws.Range("B19").Value = "Other Assets"
$ws.Range("C19").Value = "Tangible Other Assets"
$ws.Range("E19").Value = 1860
$ws.Range("B20").Value = "Total Assets"
$ws.Range("C20").Value = ""
$ws.Range("E20").Value = 129190
$ws.Range("A21").Value = "Liabilities & Shareholders' Equity"
$ws.Range("B21").Value = "ST Debt & Current Portion LT Debt"
$ws.Range("C21").Value = "Short Term Debt"
$ws.Range("E21").Value = 2000
$ws.Range("A22").Value = ""
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = "Current Portion of Long Term Debt"
$ws.Range("E22").Value = 10
$ws.Range("B23").Value = "Accounts Payable"
$ws.Range("E23").Value = 436
$ws.Range("B24").Value = "Income Tax Payable"
$ws.Range("C24").Value = ""
$ws.Range("E24").Value = 96
$ws.Range("B25").Value = "Other Current Liabilities"
$ws.Range("C25").Value = "Accrued Payroll"
$ws.Range("E25").Value = 3070
$ws.Range("B26").Value = ""
$ws.Range("C26").Value = "Miscellaneous Current Liabilities"
$ws.Range("E26").Value = 9890
$ws.Range("B27").Value = "Total Current Liabilities"
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = 16780
$ws.Range("B28").Value = "Long-Term Debt"
$ws.Range("C28").Value = "Long-Term Debt excl. Capitalized Leases"
$ws.Range("D28").Value = "Non-Convertible Debt"
$ws.Range("E28").Value = 2990
$ws.Range("B29").Value = ""
$ws.Range("C29").Value = "Capitalized Lease Obligations"
$ws.Range("E29").Value = 236
$ws.Range("B30").Value = "Deferred Taxes"
$ws.Range("C30").Value = "Deferred Taxes - Credit"
$ws.Range("E30").Value = 758
$ws.Range("B31").Value = ""
$ws.Range("C31").Value = "Deferred Taxes - Debit"
$ws.Range("E31").Value = 176
$ws.Range("B32").Value = "Other Liabilities"
$ws.Range("C32").Value = "Other Liabilities (excl. Deferred Income)"
$ws.Range("E32").Value = 4460
$ws.Range("B33").Value = ""
$ws.Range("C33").Value = "Deferred Income"
$ws.Range("E33").Value = 104
$ws.Range("B34").Value = "Total Liabilities"
$ws.Range("C34").Value = ""
$ws.Range("E34").Value = 25330
$ws.Range("B35").Value = "Common Equity (Total)"
$ws.Range("C35").Value = "Common Stock Par/Carry Value"
$ws.Range("E35").Value = 0.68
$ws.Range("C36").Value = "Additional Paid-In Capital/Capital Surplus"
$ws.Range("E36").Value = ""
$ws.Range("C37").Value = "Retained Earnings"
$ws.Range("E37").Value = 75070
$ws.Range("C38").Value = "Cumulative Translation Adjustment/Unrealized For. Exch. Gain"
$ws.Range("E38").Value = -980
$ws.Range("C39").Value = "Unrealized Gain/Loss Marketable Securities"
$ws.Range("E39").Value = 421
$ws.Range("B40").Value = ""
$ws.Range("C40").Value = "Other Appropriated Reserves"
$ws.Range("E40").Value = ""
$ws.Range("B41").Value = "Total Shareholders' Equity"
$ws.Range("E41").Value = 103860
$ws.Range("B42").Value = "Total Equity"
$ws.Range("E42").Value = 103860
$ws.Range("A43").Value = ""
$ws.Range("B43").Value = "Liabilities & Shareholders' Equity"
$ws.Range("C43").Value = ""
$ws.Range("D43").Value = ""
$ws.Range("E43").Value = 129190
